$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the title text in A1 (shortened wording)
$ws.Range("A1").Value = 'Nomenclatura para el "drone_number" '

# The URI list shifted up by one: E7E7E7E7E0 moves from B10 up to B3,
# and E1..E7 each shift down one row; D0..D3 stay put.
$ws.Range("B3").Value  = "radio://0/80/2M/E7E7E7E7E0"
$ws.Range("B4").Value  = "radio://0/80/2M/E7E7E7E7E1"
$ws.Range("B5").Value  = "radio://0/80/2M/E7E7E7E7E2"
$ws.Range("B6").Value  = "radio://0/80/2M/E7E7E7E7E3"
$ws.Range("B7").Value  = "radio://0/80/2M/E7E7E7E7E4"
$ws.Range("B8").Value  = "radio://0/80/2M/E7E7E7E7E5"
$ws.Range("B9").Value  = "radio://0/80/2M/E7E7E7E7E6"
$ws.Range("B10").Value = "radio://0/80/2M/E7E7E7E7E7"
$ws.Range("B11").Value = "radio://0/80/2M/E7E7E7E7D0"
$ws.Range("B12").Value = "radio://0/80/2M/E7E7E7E7D1"
$ws.Range("B13").Value = "radio://0/80/2M/E7E7E7E7D2"
$ws.Range("B14").Value = "radio://0/80/2M/E7E7E7E7D3"

# Move the active selection to A2
$ws.Range("A2").Select()
